{"js": "// Remove the trailing \"Ver no Jupiter ...\" / copyright footer block from the\n// document, along with the blank paragraph that separated it from the\n// \"Requisitos\" section above. These were the last three non-empty/blank\n// paragraphs before the final blank paragraph + page-break paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find the paragraph that holds the \"Ver no Jupiter ...\" text; the blank\n// paragraph immediately preceding it and the copyright paragraph\n// immediately following it are removed together with it.\nlet jupiterIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targets[0]) {\n    jupiterIndex = i;\n    break;\n  }\n}\n\nif (jupiterIndex !== -1) {\n  const copyrightIndex = jupiterIndex + 1;\n  const blankIndex = jupiterIndex - 1;\n\n  // Delete from the highest index down so earlier indices stay valid.\n  if (\n    copyrightIndex < paragraphs.items.length &&\n    paragraphs.items[copyrightIndex].text === targets[1]\n  ) {\n    paragraphs.items[copyrightIndex].delete();\n  }\n  paragraphs.items[jupiterIndex].delete();\n  if (blankIndex >= 0 && paragraphs.items[blankIndex].text === \"\") {\n    paragraphs.items[blankIndex].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / copyright footer block from the\n# document, along with the blank paragraph that separated it from the\n# \"Requisitos\" section above.\n$d = $word.ActiveDocument\n\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$jupiterIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n    if ($t -eq $jupiterText) {\n        $jupiterIndex = $i\n        break\n    }\n}\n\nif ($jupiterIndex -ne -1) {\n    $copyrightIndex = $jupiterIndex + 1\n    $blankIndex = $jupiterIndex - 1\n\n    # Delete from the highest paragraph index down so earlier indices stay valid.\n    if ($copyrightIndex -le $d.Paragraphs.Count) {\n        $ct = $d.Paragraphs.Item($copyrightIndex).Range.Text.TrimEnd([char]13)\n        if ($ct -eq $copyrightText) {\n            $d.Paragraphs.Item($copyrightIndex).Range.Delete()\n        }\n    }\n    $d.Paragraphs.Item($jupiterIndex).Range.Delete()\n    if ($blankIndex -ge 1) {\n        $bt = $d.Paragraphs.Item($blankIndex).Range.Text.TrimEnd([char]13)\n        if ($bt -eq \"\") {\n            $d.Paragraphs.Item($blankIndex).Range.Delete()\n        }\n    }\n}\n"}
